$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the template row (259) down into the 9 new rows (260-268), then overwrite
# the per-row values (columns A and C). B, D, E, H, I, J, L, M stay identical to
# row 259's pattern except H and M which get new shared text.
$ws.Range("A259:M259").Copy() | Out-Null
for ($i = 0; $i -lt 9; $i++) {
    $r = 260 + $i
    $ws.Range("A$r`:M$r").PasteSpecial(-4104) | Out-Null
}
$excel.CutCopyMode = 0

$profileNames = @(
    "INNBYGGERPOST_DPI_DIGITAL_1_0",
    "INNBYGGERPOST_DPI_UTSKRIFT_1_0",
    "INNBYGGERPOST_DPI_FLYTTET_1_0",
    "INNBYGGERPOST_DPI_LEVERINGSKVITTERING_1_0",
    "INNBYGGERPOST_DPI_FEILKVITTERING_1_0",
    "INNBYGGERPOST_DPI_AAPNINGSKVITTERING_1_0",
    "INNBYGGERPOST_DPI_MOTTAKSKVITTERING_1_0",
    "INNBYGGERPOST_DPI_VARSLINGFEILETKVITTERING_1_0",
    "INNBYGGERPOST_DPI_RETURPOSTKVITTERING_1_0"
)

$docTypeValues = @(
    "urn:fdc:digdir.no:2020:innbyggerpost:xsd::innbyggerpost##urn:fdc:digdir.no:2020:innbyggerpost:schema:digital::1.0",
    "urn:fdc:digdir.no:2020:innbyggerpost:xsd::innbyggerpost##urn:fdc:digdir.no:2020:innbyggerpost:schema:utskrift::1.0",
    "urn:fdc:digdir.no:2020:innbyggerpost:xsd::innbyggerpost##urn:fdc:digdir.no:2020:innbyggerpost:schema:flyttet::1.0",
    "urn:fdc:digdir.no:2020:innbyggerpost:xsd::innbyggerpost##urn:fdc:digdir.no:2020:innbyggerpost:schema:leveringskvittering::1.0",
    "urn:fdc:digdir.no:2020:innbyggerpost:xsd::innbyggerpost##urn:fdc:digdir.no:2020:innbyggerpost:schema:feil::1.0",
    "urn:fdc:digdir.no:2020:innbyggerpost:xsd::innbyggerpost##urn:fdc:digdir.no:2020:innbyggerpost:schema:aapningskvittering::1.0",
    "urn:fdc:digdir.no:2020:innbyggerpost:xsd::innbyggerpost##urn:fdc:digdir.no:2020:innbyggerpost:schema:mottakskvittering::1.0",
    "urn:fdc:digdir.no:2020:innbyggerpost:xsd::innbyggerpost##urn:fdc:digdir.no:2020:innbyggerpost:schema:varslingfeiletkvittering::1.0",
    "urn:fdc:digdir.no:2020:innbyggerpost:xsd::innbyggerpost##urn:fdc:digdir.no:2020:innbyggerpost:schema:returpostkvittering::1.0"
)

for ($i = 0; $i -lt 9; $i++) {
    $r = 260 + $i
    $ws.Cells.Item($r, 1).Value = $profileNames[$i]
    $ws.Cells.Item($r, 3).Value = $docTypeValues[$i]
    $ws.Cells.Item($r, 8).Value = "TICC-331"
    $ws.Cells.Item($r, 13).Value = "cenbii-procid-ubl::fake"
}

Write-Host "done"
